# Moffitt PowerPoint template: switch the deck from 16:9 widescreen
# (12192000 x 6858000 EMU) to standard 4:3 (9144000 x 6858000 EMU) and
# rescale every hard-coded shape position/size on the slide master and
# the four slide layouts so the content keeps fitting the narrower
# slide (everything scaled by width-ratio 0.75, heights unchanged).

$p = $ppt.ActivePresentation

# EMU per point - PowerPoint COM shape geometry (Left/Top/Width/Height)
# is expressed in points, while OOXML stores EMU (1 pt = 12700 EMU).
$EMU = 12700.0

function Set-ShapeGeometryEmu {
    param(
        $Shape,
        [double]$XEmu,
        [double]$YEmu,
        [double]$CxEmu,
        [double]$CyEmu
    )
    $Shape.Left   = $XEmu  / $EMU
    $Shape.Top    = $YEmu  / $EMU
    $Shape.Width  = $CxEmu / $EMU
    $Shape.Height = $CyEmu / $EMU
}

# --- 1. Slide size: Widescreen -> Standard (4:3) -------------------------
# This flips p:sldSz to cx="9144000" cy="6858000" type="screen4x3".
$p.PageSetup.SlideSize = 1

$master = $p.SlideMaster

# --- 2. Slide Master shapes ----------------------------------------------
Set-ShapeGeometryEmu $master.Shapes.Item("Title Placeholder 1")        628650  365126  7886700 1325563
Set-ShapeGeometryEmu $master.Shapes.Item("Text Placeholder 2")         628650  1825625 7886700 4351338
Set-ShapeGeometryEmu $master.Shapes.Item("Date Placeholder 3")         628650  6356351 2057400 365125
Set-ShapeGeometryEmu $master.Shapes.Item("Footer Placeholder 4")       3028950 6356351 3086100 365125
Set-ShapeGeometryEmu $master.Shapes.Item("Slide Number Placeholder 5") 6457950 6356351 2057400 365125

# --- 3. Slide Layout 1 (Title Slide) --------------------------------------
$layout1 = $master.CustomLayouts.Item(1)
Set-ShapeGeometryEmu $layout1.Shapes.Item("Title 1")     1143000 1122363 6858000 2387600
Set-ShapeGeometryEmu $layout1.Shapes.Item("Subtitle 2")  1143000 3602038 5410200 1126270
Set-ShapeGeometryEmu $layout1.Shapes.Item("Picture 6")   7315201 6155428 1694276 560132

# --- 4. Slide Layout 2 (Title and Content) --------------------------------
$layout2 = $master.CustomLayouts.Item(2)
Set-ShapeGeometryEmu $layout2.Shapes.Item("Content Placeholder 2") 342900 1825625 8490439 4351338
Set-ShapeGeometryEmu $layout2.Shapes.Item("Picture 6")             7784124 0       1356213 1769706
Set-ShapeGeometryEmu $layout2.Shapes.Item("Straight Connector 7")  0       1779945 9144000 0
Set-ShapeGeometryEmu $layout2.Shapes.Item("Picture 8")             342900  6264033 1139190 376619

# --- 5. Slide Layout 3 (Section Header) -----------------------------------
$layout3 = $master.CustomLayouts.Item(3)
Set-ShapeGeometryEmu $layout3.Shapes.Item("Title 1")               409209 1779946 8476883 2782530
Set-ShapeGeometryEmu $layout3.Shapes.Item("Text Placeholder 2")    404446 4589464 8481645 1500187
Set-ShapeGeometryEmu $layout3.Shapes.Item("Picture 6")             7784124 0       1356213 1769706
Set-ShapeGeometryEmu $layout3.Shapes.Item("Straight Connector 7")  0       1779945 9144000 0
Set-ShapeGeometryEmu $layout3.Shapes.Item("Picture 8")             342900  6264033 1139190 376619

# --- 6. Slide Layout 4 (Two Content) --------------------------------------
$layout4 = $master.CustomLayouts.Item(4)
Set-ShapeGeometryEmu $layout4.Shapes.Item("Content Placeholder 2") 342900  1825625 4171950 4351338
Set-ShapeGeometryEmu $layout4.Shapes.Item("Content Placeholder 3") 4629150 1825625 4227635 4351338
Set-ShapeGeometryEmu $layout4.Shapes.Item("Straight Connector 7")  0       1779945 9144000 0
Set-ShapeGeometryEmu $layout4.Shapes.Item("Picture 9")             7784124 0       1356213 1769706
Set-ShapeGeometryEmu $layout4.Shapes.Item("Picture 10")            342900  6264033 1139190 376619
